# Edit script: insert 5 new price-report rows for "Vega Monumental Concepción - Nectarín"
# The workbook has a single worksheet with tabular data in A:T.
# New rows are inserted at row 353 (pushing existing rows 353-428 down to 358-433),
# and are populated with new variety/price information.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows starting at row 353 (existing rows 353..428 shift down to 358..433)
$ws.Rows("353:357").Insert()

# Common (repeated) column values for this market/region/product block
$marketId   = 11
$market     = "Vega Monumental Concepción"
$region     = "Bíobío"
$codreg     = 8
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria  = "Nectarín"
$origen     = "Región de O'Higgins"

# Data for the 5 new rows: date, variedad, calidad, volumen, min, max, promedio, unidad, precioKg, kgUnidad
$newRows = @(
    @{ Row=353; Fecha=44946; Variedad="Sun Rise"; Calidad="Especial"; Volumen=50;  Min=15000; Max=15000; Prom=15000; Unidad="`$/caja 15 kilos empedrada"; PrecioKg=1000; KgUnidad=15 },
    @{ Row=354; Fecha=44946; Variedad="Sun Rise"; Calidad="Primera";  Volumen=50;  Min=13000; Max=13000; Prom=13000; Unidad="`$/caja 15 kilos empedrada"; PrecioKg=867;  KgUnidad=15 },
    @{ Row=355; Fecha=44946; Variedad="Venus";    Calidad="Especial"; Volumen=50;  Min=14000; Max=14000; Prom=14000; Unidad="`$/caja 15 kilos empedrada"; PrecioKg=933;  KgUnidad=15 },
    @{ Row=356; Fecha=44946; Variedad="Venus";    Calidad="Primera";  Volumen=50;  Min=12000; Max=12000; Prom=12000; Unidad="`$/caja 15 kilos empedrada"; PrecioKg=800;  KgUnidad=15 },
    @{ Row=357; Fecha=44946; Variedad="Venus";    Calidad="Segunda";  Volumen=50;  Min=10000; Max=10000; Prom=10000; Unidad="`$/caja 15 kilos empedrada"; PrecioKg=667;  KgUnidad=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $marketId
    $ws.Cells.Item($row, 2).Value2  = $market
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.Fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $r.Variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.Min
    $ws.Cells.Item($row, 15).Value2 = $r.Max
    $ws.Cells.Item($row, 16).Value2 = $r.Prom
    $ws.Cells.Item($row, 17).Value2 = $r.Unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value2 = $r.KgUnidad
}

Write-Host "Insert + populate complete"
